# Update the logbook "Ideal - Remaining efforts" values for rows 21-24
# (P column becomes a hard-coded value instead of the shared IF() formula,
# Q column keeps its formula but now evaluates against the new P value).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P21").Value = 2
$ws.Range("P22").Value = 2
$ws.Range("P23").Value = 2
$ws.Range("P24").Value = 2

# Row 27 (Actual - Remaining efforts) and the chart series that reads
# Sheet1!$E$27:$Q$27 recalc automatically from the SUM formulas already
# in place, picking up the lowered P/Q totals.

# Move the sheet's active selection to match the author's final cursor
# position (also drops the stale scroll-pinned topLeftCell).
[void]$ws.Range("P25").Select()
